# Update cryptos list data (price + 1h volume change, plus a couple of
# re-ranked rows) per commit "Updated cryptos list on Wed Nov 15 18:33:11
# UTC 2023 with GitHub Actions".
#
# All Price (column D) and Volume(1h) (column E) cells on this sheet are
# stored as literal text (e.g. "37.273.88", "  +3.50%  "), not numbers.
# Excel's normal type-inference would happily reinterpret plain-looking
# numeric strings (like "0.999") as numbers when assigned via .Value, so
# we force the cell format to Text ("@") first for every Price cell we
# touch to guarantee the stored type matches the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "37.122.05"
$ws.Range("E2").Value = "  +3.15%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.030.70"
$ws.Range("E3").Value = "  +0.77%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.42%  "

# Row 5 - BNB
Set-TextValue "D5" "252.61"
$ws.Range("E5").Value = "  +4.33%  "

# Row 6 - XRP
Set-TextValue "D6" "0.642"
$ws.Range("E6").Value = "  -1.36%  "

# Row 7 - Solana
Set-TextValue "D7" "62.28"
$ws.Range("E7").Value = "  +14.73%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - OKB
Set-TextValue "D9" "58.96"
$ws.Range("E9").Value = "  +0.57%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.373"
$ws.Range("E10").Value = "  +3.95%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0747"
$ws.Range("E11").Value = "  +1.94%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -1.05%  "

# Row 13 - Polygon
Set-TextValue "D13" "0.897"
$ws.Range("E13").Value = "  +1.02%  "

# Row 14 - Chainlink
Set-TextValue "D14" "14.99"
$ws.Range("E14").Value = "  +6.71%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.328.96"
$ws.Range("E15").Value = "  +0.69%  "

# Row 16 - was Avalanche, now Polkadot
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D16" "5.48"
$ws.Range("E16").Value = "  +4.60%  "

# Row 17 - was Polkadot, now Avalanche
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D17" "20.32"
$ws.Range("E17").Value = "  +18.81%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.041.95"
$ws.Range("E18").Value = "  +0.38%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "36.980.17"
$ws.Range("E19").Value = "  +3.07%  "

# Row 20 - Litecoin
Set-TextValue "D20" "72.57"
$ws.Range("E20").Value = "  +2.48%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0868"
$ws.Range("E21").Value = "  +2.95%  "

# Row 22 - Uniswap
Set-TextValue "D22" "5.30"
$ws.Range("E22").Value = "  +3.67%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "235.31"
$ws.Range("E23").Value = "  -0.13%  "

# Row 24 - PancakeSwap
Set-TextValue "D24" "2.81"
$ws.Range("E24").Value = "  +24.59%  "

# Row 25 - Dai
Set-TextValue "D25" "0.996"
$ws.Range("E25").Value = "  -0.62%  "

# Row 26 - Toncoin
Set-TextValue "D26" "2.33"
$ws.Range("E26").Value = "  -0.33%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.52"
$ws.Range("E27").Value = "  +4.14%  "

# Row 28 - Monero
Set-TextValue "D28" "164.68"
$ws.Range("E28").Value = "  +1.09%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "19.68"
$ws.Range("E29").Value = "  +0.12%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  +1.00%  "

# Row 31 - Filecoin
Set-TextValue "D31" "5.13"
$ws.Range("E31").Value = "  +5.65%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +4.77%  "

# Row 33 - Kaspa
Set-TextValue "D33" "0.111"
$ws.Range("E33").Value = "  +23.30%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "4.69"
$ws.Range("E34").Value = "  +9.54%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.0611"
$ws.Range("E35").Value = "  +3.40%  "

# Row 36 - LidoDAOToken
Set-TextValue "D36" "2.46"
$ws.Range("E36").Value = "  +13.68%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  -0.16%  "

# Row 38 - was THORChain, now WEMIXToken
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D38" "1.81"
$ws.Range("E38").Value = "  -0.83%  "

# Row 39 - was WEMIXToken, now THORChain
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D39" "5.97"
$ws.Range("E39").Value = "  +22.43%  "

# Row 40 - Cronos
$ws.Range("E40").Value = "  +17.78%  "

# Row 41 - TrustWalletToken
Set-TextValue "D41" "1.23"
$ws.Range("E41").Value = "  +3.21%  "

# Row 42 - RenderToken
Set-TextValue "D42" "2.77"
$ws.Range("E42").Value = "  +24.01%  "

# Row 43 - HuobiToken
Set-TextValue "D43" "2.93"
$ws.Range("E43").Value = "  +1.63%  "

# Row 44 - ARBITRUM
Set-TextValue "D44" "1.13"
$ws.Range("E44").Value = "  +3.95%  "

# Row 45 - VeChain
Set-TextValue "D45" "0.0217"
$ws.Range("E45").Value = "  +2.13%  "

# Row 46 - FraxShare
Set-TextValue "D46" "8.08"
$ws.Range("E46").Value = "  +10.25%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "16.77"
$ws.Range("E47").Value = "  +9.56%  "

# Row 48 - Aave
Set-TextValue "D48" "94.86"
$ws.Range("E48").Value = "  +4.08%  "

# Row 49 - Maker
Set-TextValue "D49" "1.428.34"
$ws.Range("E49").Value = "  +3.30%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +0.97%  "

# Row 51 - MultiversX
Set-TextValue "D51" "47.23"
$ws.Range("E51").Value = "  +4.15%  "
